$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A1").Value = "test"
